$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data. Values are entered with a leading
# apostrophe to force text interpretation (several look like numbers, e.g.
# "582.28" or multi-dot "68.104.67"), then the style is reset to "Normal" so
# the cell keeps the original (unstyled) text formatting instead of picking up
# an implicit @ (Text) number format.

$ws.Range("D2").Value = "'68.104.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.251.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'582.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'184.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.98%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.30%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.418"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.828.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.216.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.47%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.42%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.261.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'5.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.39%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'392.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'Dai"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.20%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'Litecoin"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'71.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.519"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.84%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +4.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'PancakeSwap"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.62%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'5.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'22.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.23%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.64%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'163.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.71%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +5.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.823"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.65%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'26.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.55%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.20%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.75%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.54%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'25.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.82%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0689"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'41.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.33%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.652.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.61%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'337.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0282"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.84%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'31.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.30%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.44%  "
$ws.Range("E51").Style = "Normal"
